# Fruta / hortaliza, semanal
# Insert two new weekly price rows at the top of the "Vega Monumental
# Concepción - Naranja" data block (old rows 199-226 shift down to 201-228),
# then populate the two newly inserted rows (199-200) with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 199-226 down by two rows.
$ws.Rows.Item(199).Insert()
$ws.Rows.Item(199).Insert()

# --- New row 199: Naranja, Valencia, Primera ---
$ws.Cells.Item(199, 1).Value = 11
$ws.Cells.Item(199, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(199, 3).Value = "Bíobío"
$ws.Cells.Item(199, 4).Value = 44617
$ws.Cells.Item(199, 5).Value = 8
$ws.Cells.Item(199, 6).Value = "Fruta"
$ws.Cells.Item(199, 7).Value = 100102
$ws.Cells.Item(199, 8).Value = "Cítricos"
$ws.Cells.Item(199, 9).Value = 100102005
$ws.Cells.Item(199, 10).Value = "Naranja"
$ws.Cells.Item(199, 11).Value = "Valencia"
$ws.Cells.Item(199, 12).Value = "Primera"
$ws.Cells.Item(199, 13).Value = 300
$ws.Cells.Item(199, 14).Value = 9000
$ws.Cells.Item(199, 15).Value = 10000
$ws.Cells.Item(199, 16).Value = 9667
$ws.Cells.Item(199, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(199, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(199, 19).Value = 644
$ws.Cells.Item(199, 20).Value = 15

# --- New row 200: Naranja, Valencia, Segunda ---
$ws.Cells.Item(200, 1).Value = 11
$ws.Cells.Item(200, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(200, 3).Value = "Bíobío"
$ws.Cells.Item(200, 4).Value = 44617
$ws.Cells.Item(200, 5).Value = 8
$ws.Cells.Item(200, 6).Value = "Fruta"
$ws.Cells.Item(200, 7).Value = 100102
$ws.Cells.Item(200, 8).Value = "Cítricos"
$ws.Cells.Item(200, 9).Value = 100102005
$ws.Cells.Item(200, 10).Value = "Naranja"
$ws.Cells.Item(200, 11).Value = "Valencia"
$ws.Cells.Item(200, 12).Value = "Segunda"
$ws.Cells.Item(200, 13).Value = 100
$ws.Cells.Item(200, 14).Value = 8000
$ws.Cells.Item(200, 15).Value = 8000
$ws.Cells.Item(200, 16).Value = 8000
$ws.Cells.Item(200, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(200, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(200, 19).Value = 533
$ws.Cells.Item(200, 20).Value = 15
